$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1476.8667
$ws.Range("J19").Value = 1732.909
$ws.Range("L19").Value = 1732.909
$ws.Range("N19").Value = -2082.909

# Row 41
$ws.Range("H41").Value = 244.23529
$ws.Range("I41").Value = 257.4
$ws.Range("J41").Value = 225.42857
$ws.Range("K41").Value = 257.4
$ws.Range("L41").Value = 225.42857
$ws.Range("M41").Value = 182.6
$ws.Range("N41").Value = -1105.42857

# Row 107
$ws.Range("H107").Value = 1001.45
$ws.Range("I107").Value = 1016.64703
$ws.Range("J107").Value = 915.3333
$ws.Range("K107").Value = 1016.64703
$ws.Range("L107").Value = 915.3333
$ws.Range("M107").Value = 903.35297
$ws.Range("N107").Value = -4755.3333

# Row 118
$ws.Range("H118").Value = 1049.5
$ws.Range("I118").Value = 1066
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 3198
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = -1541
$ws.Range("N118").Value = -6314

# Row 125
$ws.Range("H125").Value = 1696.5555
$ws.Range("I125").Value = 1074.75
$ws.Range("J125").Value = 2194
$ws.Range("K125").Value = 9672.75
$ws.Range("L125").Value = 19746
$ws.Range("M125").Value = -7212.75
$ws.Range("N125").Value = -24666

# Row 127
$ws.Range("H127").Value = 1021.44446
$ws.Range("I127").Value = 598.8333
$ws.Range("J127").Value = 1866.6666
$ws.Range("K127").Value = 1796.4999
$ws.Range("L127").Value = 5599.9998
$ws.Range("M127").Value = 3163.5001
$ws.Range("N127").Value = -15519.9998

# Row 137
$ws.Range("H137").Value = 31532950
$ws.Range("I137").Value = 6667981.5
$ws.Range("J137").Value = 83334970
$ws.Range("K137").Value = 20003944.5
$ws.Range("L137").Value = 250004910
$ws.Range("M137").Value = -20001394.5
$ws.Range("N137").Value = -250010010

# Row 138
$ws.Range("H138").Value = 2418.0378
$ws.Range("I138").Value = 1984.4333
$ws.Range("J138").Value = 2983.6086
$ws.Range("K138").Value = 5953.2999
$ws.Range("L138").Value = 8950.825800000001
$ws.Range("M138").Value = -813.2999
$ws.Range("N138").Value = -19230.8258

# Row 141
$ws.Range("H141").Value = 1144.2667
$ws.Range("I141").Value = 732.2162
$ws.Range("K141").Value = 2196.6486
$ws.Range("M141").Value = 2983.3514


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 10841990
$ws.Range("I61").Value = 10841990
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10841990
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -10841778
$ws.Range("N61").ClearContents()

# Row 74
$ws.Range("H74").Value = 35638110
$ws.Range("I74").Value = 60785296
$ws.Range("J74").Value = 12936.917
$ws.Range("K74").Value = 60785296
$ws.Range("L74").Value = 12936.917
$ws.Range("M74").Value = -60784422
$ws.Range("N74").Value = -14684.917

# Row 77
$ws.Range("H77").Value = 35638110
$ws.Range("I77").Value = 60785296
$ws.Range("J77").Value = 12936.917
$ws.Range("K77").Value = 303926480
$ws.Range("L77").Value = 64684.585
$ws.Range("M77").Value = -303922112
$ws.Range("N77").Value = -73420.58499999999

# Row 102
$ws.Range("H102").Value = 1999.4736
$ws.Range("I102").Value = 1994
$ws.Range("J102").Value = 2020
$ws.Range("K102").Value = 1994
$ws.Range("L102").Value = 2020
$ws.Range("M102").Value = -372
$ws.Range("N102").Value = -5264

# Row 122
$ws.Range("H122").Value = 1287.2
$ws.Range("I122").Value = 1244.0952
$ws.Range("J122").Value = 1387.7778
$ws.Range("K122").Value = 3732.2856
$ws.Range("L122").Value = 4163.3334
$ws.Range("M122").Value = -1282.2856
$ws.Range("N122").Value = -9063.3334

# Row 132
$ws.Range("H132").Value = 1185043.6
$ws.Range("I132").Value = 1644537.4
$ws.Range("J132").Value = 112891.78
$ws.Range("K132").Value = 4933612.199999999
$ws.Range("L132").Value = 338675.34
$ws.Range("M132").Value = -4931082.199999999
$ws.Range("N132").Value = -343735.34

# Row 136
$ws.Range("H136").Value = 10841990
$ws.Range("I136").Value = 10841990
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 32525970
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -32523420
$ws.Range("N136").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 460425.12
$ws.Range("I31").Value = 1221.92
$ws.Range("J31").Value = 1280430.9
$ws.Range("K31").Value = 1221.92
$ws.Range("L31").Value = 1280430.9
$ws.Range("M31").Value = -926.9200000000001
$ws.Range("N31").Value = -1281020.9

# Row 34
$ws.Range("H34").Value = 460425.12
$ws.Range("I34").Value = 1221.92
$ws.Range("J34").Value = 1280430.9
$ws.Range("K34").Value = 1221.92
$ws.Range("L34").Value = 1280430.9
$ws.Range("M34").Value = -1019.92
$ws.Range("N34").Value = -1280834.9

# Row 58
$ws.Range("H58").Value = 1738.8441
$ws.Range("I58").Value = 786
$ws.Range("J58").Value = 3843.0417
$ws.Range("K58").Value = 786
$ws.Range("L58").Value = 3843.0417
$ws.Range("M58").Value = -583
$ws.Range("N58").Value = -4249.0417

# Row 132
$ws.Range("H132").Value = 1661.16
$ws.Range("I132").Value = 1697.3784
$ws.Range("J132").Value = 1558.0769
$ws.Range("K132").Value = 5092.135200000001
$ws.Range("L132").Value = 4674.2307
$ws.Range("M132").Value = -2562.135200000001
$ws.Range("N132").Value = -9734.2307

# Row 136
$ws.Range("H136").Value = 1738.8441
$ws.Range("I136").Value = 786
$ws.Range("J136").Value = 3843.0417
$ws.Range("K136").Value = 2358
$ws.Range("L136").Value = 11529.1251
$ws.Range("M136").Value = 192
$ws.Range("N136").Value = -16629.1251


$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 76.26316
$ws.Range("I40").Value = 76.26316
$ws.Range("K40").Value = 305.05264
$ws.Range("M40").Value = -236.05264

# Row 106
$ws.Range("H106").Value = 5517.391
$ws.Range("J106").Value = 5517.391
$ws.Range("L106").Value = 16552.173
$ws.Range("N106").Value = -18444.173

# Row 121
$ws.Range("H121").Value = 25000578
$ws.Range("I121").Value = 318.33334
$ws.Range("J121").Value = 62500970
$ws.Range("K121").Value = 955.0000200000001
$ws.Range("L121").Value = 187502910
$ws.Range("M121").Value = 354.9999799999999
$ws.Range("N121").Value = -187505530

# Row 131
$ws.Range("H131").Value = 2476.5908
$ws.Range("I131").Value = 7078.8887
$ws.Range("J131").Value = 1749.9122
$ws.Range("K131").Value = 21236.6661
$ws.Range("L131").Value = 5249.7366
$ws.Range("M131").Value = -16196.6661
$ws.Range("N131").Value = -15329.7366

# Row 134
$ws.Range("H134").Value = 2722.2693
$ws.Range("I134").Value = 2425.1738
$ws.Range("K134").Value = 7275.5214
$ws.Range("M134").Value = -2205.5214

# Row 137
$ws.Range("H137").Value = 19824.213
$ws.Range("I137").Value = 2122.353
$ws.Range("J137").Value = 26663.568
$ws.Range("K137").Value = 6367.059
$ws.Range("L137").Value = 79990.704
$ws.Range("M137").Value = -1267.059
$ws.Range("N137").Value = -90190.704


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4471.3335
$ws.Range("I122").Value = 3723.5
$ws.Range("K122").Value = 11170.5
$ws.Range("M122").Value = -8720.5

# Row 132
$ws.Range("H132").Value = 1815.1714
$ws.Range("I132").Value = 1853
$ws.Range("J132").Value = 1687.5
$ws.Range("K132").Value = 5559
$ws.Range("L132").Value = 5062.5
$ws.Range("M132").Value = -3029
$ws.Range("N132").Value = -10122.5


$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 5325.7334
$ws.Range("I136").Value = 5454.6665
$ws.Range("J136").Value = 4810
$ws.Range("K136").Value = 16363.9995
$ws.Range("L136").Value = 14430
$ws.Range("M136").Value = -13813.9995
$ws.Range("N136").Value = -19530


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2327.76
$ws.Range("I122").Value = 1344.1111
$ws.Range("J122").Value = 4857.143
$ws.Range("K122").Value = 4032.3333
$ws.Range("L122").Value = 14571.429
$ws.Range("M122").Value = -1582.3333
$ws.Range("N122").Value = -19471.429

# Row 138
$ws.Range("H138").Value = 49800
$ws.Range("J138").Value = 49800
$ws.Range("L138").Value = 49800
$ws.Range("N138").Value = -60080

